$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New "ActualRate" (column E) values for rows 2-31. These now equal the
# "ExpectedRate" (column D) values for each row, formatted as currency text,
# and the "Result" (column F) becomes "PASS" for every row (New Pre-Prod URL).
$actualRates = @{
    2  = "$19.04"
    3  = "$24.32"
    4  = "$28.55"
    5  = "$40.19"
    6  = "$52.88"
    7  = "$65.57"
    8  = "$77.91"
    9  = "$231.91"
    10 = "$15.86"
    11 = "$20.09"
    12 = "$24.32"
    13 = "$37.01"
    14 = "$49.70"
    15 = "$65.57"
    16 = "$12.69"
    17 = "$15.86"
    18 = "$19.04"
    19 = "$29.61"
    20 = "$40.19"
    21 = "$48.65"
    22 = "$38.07"
    23 = "$48.65"
    24 = "$57.11"
    25 = "$80.37"
    26 = "$105.75"
    27 = "$105.75"
    28 = "$150.17"
    29 = "$317.25"
    30 = "$473.23"
    31 = "$252.98"
}

foreach ($row in 2..31) {
    $eCell = $ws.Cells.Item($row, 5)
    # Force the value to be written as literal text (not auto-converted to a
    # currency number by Excel's input parser), then restore the cell's
    # original (default) style so no stray style id is introduced.
    $eCell.NumberFormat = "@"
    $eCell.Value = $actualRates[$row]
    $eCell.Style = "Normal"

    $ws.Cells.Item($row, 6).Value = "PASS"
}
